# Weekly update: a new "Espárragos" price record for Femacal de La Calera
# (Coquimbo -> Provincia de Quillota) is inserted as row 33 of the data
# table, pushing the existing rows 33-58 down to 34-59 (dimension grows
# from A1:R58 to A1:R59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 33, shifting everything below
# (including formatting) down by one - this is what lets row 34's date
# cell keep its existing "s=2" (date) style, etc.
$ws.Rows.Item(33).Insert()

# Populate the new row 33 with the new weekly record.
$ws.Range('A33').Value = 3
$ws.Range('B33').Value = 'Femacal de La Calera'
$ws.Range('C33').Value = 'Coquimbo'
$ws.Range('D33').Value = 44904
$ws.Range('E33').Value = 5
$ws.Range('F33').Value = 300000000
$ws.Range('G33').Value = 'Espárragos'
$ws.Range('H33').Value = 'Verde'
$ws.Range('I33').Value = 'Primera'
$ws.Range('J33').Value = 1100
$ws.Range('K33').Value = 1400
$ws.Range('L33').Value = 1400
$ws.Range('M33').Value = 1400
$ws.Range('N33').Value = '$/kilo'
$ws.Range('O33').Value = 'Provincia de Quillota'
$ws.Range('P33').Value = 1400
$ws.Range('Q33').Value = 1
$ws.Range('R33').Value = 'Hortaliza'

# Make sure the new date cell carries the same number format the other
# "Fecha" cells in column D use.
$ws.Range('D33').NumberFormat = $ws.Range('D34').NumberFormat
